$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Delete the three now-unneeded expense rows (old "Egreso 3", "Egreso 4", "Egreso 4") ---
# This shifts old row 14 (Total Egresos) -> row 11 and old row 15 (Acumulado) -> row 12.
$ws.Rows.Item(11).Resize(3).Delete()

# --- Add new column E (Mes 4) with the same width as column D ---
$ws.Columns.Item(5).ColumnWidth = $ws.Columns.Item(4).ColumnWidth

# --- Copy the formatting (fill + border) of the styled column-D cells onto the new column-E cells ---
$ws.Range("D7").Copy()
$ws.Range("E7").PasteSpecial(-4122)
$ws.Range("D11").Copy()
$ws.Range("E11").PasteSpecial(-4122)
$ws.Range("D12").Copy()
$ws.Range("E12").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Header row ---
$ws.Range("E1").Value = "Mes 4"

# --- Ingresos section ---
# Row 3: Prestamo
$ws.Range("B3").Value = 12500
$ws.Range("C3").Value = 0
$ws.Range("D3").Value = 0
$ws.Range("E3").Value = 0

# Row 4: Donaciones
$ws.Range("B4").Value = 1200000
$ws.Range("C4").Value = 0
$ws.Range("D4").Value = 0
$ws.Range("E4").Value = 0

# Row 5: Patrocinador
$ws.Range("B5").Value = 0
$ws.Range("C5").Value = 0
$ws.Range("D5").Value = 0
$ws.Range("E5").Value = 0

# Row 6: Pago Cliente
$ws.Range("B6").Value = 0
$ws.Range("C6").Value = 0
$ws.Range("D6").Value = 0
$ws.Range("E6").Value = 0

# Row 7: Total Ingresos (styled row)
$ws.Range("B7").Value = 1212500
$ws.Range("C7").Value = 0
$ws.Range("D7").Value = 0
$ws.Range("E7").Value = 0

# --- Egresos section ---
# Row 9: renamed "Egreso 1" -> "Ingenieros"
$ws.Range("A9").Value = "Ingenieros"
$ws.Range("B9").Value = 2400
$ws.Range("C9").Value = 0
$ws.Range("D9").Value = 0
$ws.Range("E9").Value = 0

# Row 10: renamed "Egreso 2" -> "Ingenieros"
$ws.Range("A10").Value = "Ingenieros"
$ws.Range("B10").Value = 0
$ws.Range("C10").Value = 3600
$ws.Range("D10").Value = 0
$ws.Range("E10").Value = 0

# Row 11: Total Egresos (styled row, shifted up after deletion)
$ws.Range("B11").Value = 2400
$ws.Range("C11").Value = 3600
$ws.Range("D11").Value = 0
$ws.Range("E11").Value = 0

# Row 12: Acumulado (styled row, shifted up after deletion)
$ws.Range("B12").Value = 1210100
$ws.Range("C12").Value = -3600
$ws.Range("D12").Value = 0
$ws.Range("E12").Value = 0
